# Update "想去人数" (want-to-go count) figures after a fresh data pull.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - rows 3,4,5 hold the affected events in column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 2086
$wsExpo.Range("F4").Value = 862
$wsExpo.Range("F5").Value = 1250

# Sheet "全部类型" (All types) - same events appear at rows 3, 6, 7
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 2086
$wsAll.Range("F6").Value = 862
$wsAll.Range("F7").Value = 1250
